$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = "rise"
$ws.Cells.Item(2, 4).Value = 0.7669572830200195
$ws.Cells.Item(3, 3).Value = "demand"
$ws.Cells.Item(3, 4).Value = 0.7097986340522766
$ws.Cells.Item(4, 3).Value = "concern"
$ws.Cells.Item(4, 4).Value = 0.6822013854980469
$ws.Cells.Item(5, 3).Value = "mortgage"
$ws.Cells.Item(5, 4).Value = 0.6736873388290405
$ws.Cells.Item(6, 3).Value = "high"
$ws.Cells.Item(6, 4).Value = 0.673154890537262
$ws.Cells.Item(7, 3).Value = "potential"
$ws.Cells.Item(7, 4).Value = 0.6590673923492432
$ws.Cells.Item(8, 3).Value = "inventory"
$ws.Cells.Item(8, 4).Value = 0.6569727063179016
$ws.Cells.Item(9, 3).Value = "although"
$ws.Cells.Item(9, 4).Value = 0.6566945314407349
$ws.Cells.Item(10, 3).Value = "despite"
$ws.Cells.Item(10, 4).Value = 0.6519894003868103
$ws.Cells.Item(11, 3).Value = "because"
$ws.Cells.Item(11, 4).Value = 0.6513103842735291
$ws.Cells.Item(12, 3).Value = "employment"
$ws.Cells.Item(12, 4).Value = 0.5911222696304321
$ws.Cells.Item(13, 3).Value = "pressure"
$ws.Cells.Item(13, 4).Value = 0.5860975384712219
$ws.Cells.Item(14, 3).Value = "upward"
$ws.Cells.Item(14, 4).Value = 0.575556218624115
$ws.Cells.Item(15, 3).Value = "growth"
$ws.Cells.Item(15, 4).Value = 0.5642092227935791
$ws.Cells.Item(16, 3).Value = "however"
$ws.Cells.Item(16, 4).Value = 0.5639440417289734
$ws.Cells.Item(17, 3).Value = "contact"
$ws.Cells.Item(17, 4).Value = 0.5597914457321167
$ws.Cells.Item(18, 3).Value = "increase"
$ws.Cells.Item(18, 4).Value = 0.5471316576004028
$ws.Cells.Item(19, 3).Value = "overall"
$ws.Cells.Item(19, 4).Value = 0.5403355360031128
$ws.Cells.Item(20, 3).Value = "balance"
$ws.Cells.Item(20, 4).Value = 0.5332819819450378
$ws.Cells.Item(21, 3).Value = "continue"
$ws.Cells.Item(21, 4).Value = 0.5286637544631958
$ws.Cells.Item(22, 3).Value = "investment"
$ws.Cells.Item(22, 4).Value = 0.5767127275466919
$ws.Cells.Item(23, 3).Value = "capital"
$ws.Cells.Item(23, 4).Value = 0.534135639667511
$ws.Cells.Item(24, 3).Value = "automation"
$ws.Cells.Item(24, 4).Value = 0.5193997025489807
$ws.Cells.Item(25, 3).Value = "technology"
$ws.Cells.Item(25, 4).Value = 0.5183722376823425
$ws.Cells.Item(26, 3).Value = "expenditure"
$ws.Cells.Item(26, 4).Value = 0.5148558020591736
$ws.Cells.Item(27, 3).Value = "plan"
$ws.Cells.Item(27, 4).Value = 0.4955677688121795
$ws.Cells.Item(28, 3).Value = "develop"
$ws.Cells.Item(28, 4).Value = 0.4913558959960937
$ws.Cells.Item(29, 3).Value = "development"
$ws.Cells.Item(29, 4).Value = 0.4828211069107055
$ws.Cells.Item(30, 3).Value = "comply"
$ws.Cells.Item(30, 4).Value = 0.4719538390636444
$ws.Cells.Item(31, 3).Value = "respirator"
$ws.Cells.Item(31, 4).Value = 0.4682314097881317
$ws.Cells.Item(32, 3).Value = "skilled"
$ws.Cells.Item(32, 4).Value = 0.6720280647277832
$ws.Cells.Item(33, 3).Value = "occupation"
$ws.Cells.Item(33, 4).Value = 0.5966830253601074
$ws.Cells.Item(34, 3).Value = "worker"
$ws.Cells.Item(34, 4).Value = 0.5698993802070618
$ws.Cells.Item(35, 3).Value = "position"
$ws.Cells.Item(35, 4).Value = 0.5390143394470215
$ws.Cells.Item(36, 3).Value = "skill"
$ws.Cells.Item(36, 4).Value = 0.535441517829895
$ws.Cells.Item(37, 3).Value = "engineer"
$ws.Cells.Item(37, 4).Value = 0.5307133793830872
$ws.Cells.Item(38, 3).Value = "programmer"
$ws.Cells.Item(38, 4).Value = 0.5272024869918823
$ws.Cells.Item(39, 3).Value = "hourly"
$ws.Cells.Item(39, 4).Value = 0.5150139331817627
$ws.Cells.Item(40, 3).Value = "tradespeople"
$ws.Cells.Item(40, 4).Value = 0.5116386413574219
$ws.Cells.Item(41, 3).Value = "assistant"
$ws.Cells.Item(41, 4).Value = 0.5057680606842041
$ws.Cells.Item(42, 3).Value = "cautious"
$ws.Cells.Item(42, 4).Value = 0.6335276365280151
$ws.Cells.Item(43, 3).Value = "negative"
$ws.Cells.Item(43, 4).Value = 0.6159318089485168
$ws.Cells.Item(44, 3).Value = "pessimistic"
$ws.Cells.Item(44, 4).Value = 0.5981963872909546
$ws.Cells.Item(45, 3).Value = "optimistic"
$ws.Cells.Item(45, 4).Value = 0.5845582485198975
$ws.Cells.Item(46, 3).Value = "regard"
$ws.Cells.Item(46, 4).Value = 0.5608965754508972
$ws.Cells.Item(47, 3).Value = "uncertainty"
$ws.Cells.Item(47, 4).Value = 0.55103999376297
$ws.Cells.Item(48, 3).Value = "outlook"
$ws.Cells.Item(48, 4).Value = 0.5505427122116089
$ws.Cells.Item(49, 3).Value = "weak"
$ws.Cells.Item(49, 4).Value = 0.5421721339225769
$ws.Cells.Item(50, 3).Value = "whether"
$ws.Cells.Item(50, 4).Value = 0.5401925444602966
$ws.Cells.Item(51, 3).Value = "express"
$ws.Cells.Item(51, 4).Value = 0.5181381702423096

Write-Host "Updated similar words and similarity scores"